$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix a pre-existing inconsistency: the "inputDate" row's Parameters cell used
# "Object: String" (with a stray space) while every other row spells it
# "Object:String". Normalizing this makes the old, now-duplicate shared
# string unused so it is dropped from the shared-strings table on save.
$ws.Range("C8").Value = "Object:String"

# Populate the previously-blank row 14 with the new "clickAlert" action
# keyword. Copy formatting from row 13 (same section) first so the new row
# matches the look of the rest of the "General" keyword block, then set the
# actual cell text.
$ws.Range("A13:G13").Copy($ws.Range("A14:G14"))

$ws.Range("B14").Value = "clickAlert"
$ws.Range("C14").Value = "Object:String"
$ws.Range("E14").Value = "ObjectKey"
$ws.Range("F14").Value = "clickAlert"
$ws.Range("G14").Value = "OK"
$ws.Range("D14").Value = "Alert box actions (yes/no)"

# A14 (category column) stays blank, matching rows 3-13.
$ws.Range("A14").ClearContents()

# Move the active selection to D17, reflecting where the author's cursor
# ended up after the edit.
[void]$ws.Range("D17").Select()
